$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update shared strings for cells A50:A73: "UCD Dublin" -> "UCD Comp Sci"
$ws.Range("A50").Value = "Hi I'm a Twitterbot from UCD Comp Sci, reply ""YES"" if you would you like a story?"
$ws.Range("A51").Value = "Hi I'm a Twitterbot from UCD Comp Sci, reply ""YES"" if you want a story?"
$ws.Range("A52").Value = "Hi I'm a Twitterbot from UCD Comp Sci, reply ""YES"" if you fancy a story?"
$ws.Range("A53").Value = "Hi I'm a Twitterbot from UCD Comp Sci, reply ""YES"" if you’re up for a story?"
$ws.Range("A54").Value = "Hi I'm a Twitterbot from UCD Comp Sci, reply ""YES"" if you would like me to tell you a story?"
$ws.Range("A55").Value = "Hi I'm a Twitterbot from UCD Comp Sci, reply ""YES"" if you want me to tell you a story?"
$ws.Range("A56").Value = "Hi I'm a Twitterbot from UCD Comp Sci, reply ""YES"" if you fancy a story told by me?"
$ws.Range("A57").Value = "Hi I'm a Twitterbot from UCD Comp Sci, reply ""YES"" if you’re up for me telling you a story?"
$ws.Range("A58").Value = "Hi I'm a Twitterbot from UCD Comp Sci, reply ""YES"" if  you would like me to tweet you story?"
$ws.Range("A59").Value = "Hi I'm a Twitterbot from UCD Comp Sci, reply ""YES"" if you want a story in tweets?"
$ws.Range("A60").Value = "Hi I'm a Twitterbot from UCD Comp Sci, reply ""YES"" if you fancy a story in tweets?"
$ws.Range("A61").Value = "Hi I'm a Twitterbot from UCD Comp Sci, reply ""YES"" if you’re up for me tweeting you a story?"
$ws.Range("A62").Value = "Hi I'm a Twitterbot from UCD Comp Sci, reply ""YES"" if you would like an entertaining story?"
$ws.Range("A63").Value = "Hi I'm a Twitterbot from UCD Comp Sci, reply ""YES"" if you want an entertaining story?"
$ws.Range("A64").Value = "Hi I'm a Twitterbot from UCD Comp Sci, reply ""YES"" if you fancy an entertaining story?"
$ws.Range("A65").Value = "Hi I'm a Twitterbot from UCD Comp Sci, reply ""YES"" for an entertaining story?"
$ws.Range("A66").Value = "Hi I'm a Twitterbot from UCD Comp Sci, reply ""YES"" if you would like me to tell you an entertaining story?"
$ws.Range("A67").Value = "Hi I'm a Twitterbot from UCD Comp Sci, reply ""YES"" if you want me to tell you an entertaining story?"
$ws.Range("A68").Value = "Hi I'm a Twitterbot from UCD Comp Sci, reply ""YES"" if you fancy an entertaining story told by me?"
$ws.Range("A69").Value = "Hi I'm a Twitterbot from UCD Comp Sci, reply ""YES"" for me to tell you an entertaining story?"
$ws.Range("A70").Value = "Hi I'm a Twitterbot from UCD Comp Sci, reply ""YES"" if you would like me to tweet you an entertaining story?"
$ws.Range("A71").Value = "Hi I'm a Twitterbot from UCD Comp Sci, reply ""YES"" if you want an entertaining story in tweets?"
$ws.Range("A72").Value = "Hi I'm a Twitterbot from UCD Comp Sci, reply ""YES"" if you fancy an entertaining story in tweets?"
$ws.Range("A73").Value = "Hi I'm a Twitterbot from UCD Comp Sci, reply ""YES"" for me to tweet you an entertaining story? "

# Update the sheet view: clear the scrolled topLeftCell and move the
# selection/active cell to A18 (a single-cell selection) instead of the
# previous A65 active cell over an A1:A73 selection range.
$ws.Range("A18").Select() | Out-Null
